# "combine sample xlsx files"
#
# Rename the original sheet to "normal", then append two more sheets
# ("wrong_filename", "missing_filename") holding variants of the same
# group/filename table, and normalize the view zoom to 100%.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Sheet1 -> normal -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "normal"
$ws1.Activate()
$excel.ActiveWindow.Zoom = 100

# --- Sheet 2: wrong_filename ----------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "wrong_filename"

$data2 = @(
    @("group", "filename"),
    @("C", "image001--.png"),
    @("C", "image002--.png"),
    @("C", "image003.png"),
    @("T", "image004.png"),
    @("T", "image005.png"),
    @("T", "image006.png")
)
for ($r = 0; $r -lt $data2.Length; $r++) {
    $ws2.Cells.Item($r + 1, 1).Value = $data2[$r][0]
    $ws2.Cells.Item($r + 1, 2).Value = $data2[$r][1]
}
# column B width == 19 characters (ColumnWidth is offset by 5/MDW vs. the
# raw XML width unit, so compensate to land on an XML width of exactly 19)
$ws2.Columns.Item(2).ColumnWidth = 18.285714285714285

$ws2.Activate()
$excel.ActiveWindow.Zoom = 100

# --- Sheet 3: missing_filename ---------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "missing_filename"

$data3 = @(
    @("group", "file"),
    @("C", "image001.png"),
    @("C", "image002.png"),
    @("C", "image003.png"),
    @("T", "image004.png"),
    @("T", "image005.png"),
    @("T", "image006.png")
)
for ($r = 0; $r -lt $data3.Length; $r++) {
    $ws3.Cells.Item($r + 1, 1).Value = $data3[$r][0]
    $ws3.Cells.Item($r + 1, 2).Value = $data3[$r][1]
}
$ws3.Columns.Item(2).ColumnWidth = 18.285714285714285

$ws3.Activate()
$excel.ActiveWindow.Zoom = 100
$ws3.Range("B11").Select() | Out-Null

# Leave "normal" as the selected/active tab, matching the original workbook.
$ws1.Activate()
